# #5: cash & deposit done
#
# Rebuild the "存款" (deposit) worksheet (2nd tab) so it carries the full
# generic property schema (property_category, category, date,
# legislator_name, legislator_id, source_file, index) in addition to its
# own bank / deposit_type / currency / owner / total columns.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 汽車 - used as a safe source for the literal
                                 # "2012-02-29" text (Copy keeps it as text
                                 # instead of Excel re-parsing it as a date)
$ws2 = $wb.Worksheets.Item(2)   # 存款

# ---------------- Header row (row 1) ----------------
$ws2.Cells.Item(1, 2).Value  = "bank"
$ws2.Cells.Item(1, 3).Value  = "deposit_type"
$ws2.Cells.Item(1, 4).Value  = "currency"
$ws2.Cells.Item(1, 5).Value  = "owner"
$ws2.Cells.Item(1, 6).Value  = "total"
$ws2.Cells.Item(1, 7).Value  = "property_category"
$ws2.Cells.Item(1, 8).Value  = "category"
$ws2.Cells.Item(1, 9).Value  = "date"
$ws2.Cells.Item(1, 10).Value = "legislator_name"
$ws2.Cells.Item(1, 11).Value = "legislator_id"
$ws2.Cells.Item(1, 12).Value = "source_file"
$ws2.Cells.Item(1, 13).Value = "index"

# ---------------- Data rows (rows 2-8) ----------------
# columns: A index | B bank | C deposit_type | D currency | E owner | F total
#          G property_category | H category | I date | J legislator_name
#          K legislator_id | L source_file | M index

$rows = @(
    @{ Idx = 43; Bank = "臺灣銀行鳳山分行";         Dep = "綜合存款";     Cur = "新臺幣"; Total = 38366 },
    @{ Idx = 44; Bank = "高雄銀行市府分行";         Dep = "活期儲蓄存款"; Cur = "新臺幣"; Total = 2359237 },
    @{ Idx = 46; Bank = "中國信託商業銀行城中分行"; Dep = "活期儲蓄存款"; Cur = "新臺幣"; Total = 249260 },
    @{ Idx = 47; Bank = "中國信託商業銀行城中分行"; Dep = "活期儲蓄存款"; Cur = "美金";   Total = 729936.97 },
    @{ Idx = 48; Bank = "台新國際商業銀行鳳山分行"; Dep = "綜合存款";     Cur = "新臺幣"; Total = 3407020 },
    @{ Idx = 49; Bank = "台新國際商業銀行鳳山分行"; Dep = "綜合存款";     Cur = "歐元";   Total = 106.47 },
    @{ Idx = 50; Bank = "台新國際商業銀行鳳山分行"; Dep = "綜合存款";     Cur = "美金";   Total = 25127.62 }
)

$r = 2
foreach ($row in $rows) {
    $ws2.Cells.Item($r, 1).Value  = $row.Idx
    $ws2.Cells.Item($r, 2).Value  = $row.Bank
    $ws2.Cells.Item($r, 3).Value  = $row.Dep
    $ws2.Cells.Item($r, 4).Value  = $row.Cur
    $ws2.Cells.Item($r, 5).Value  = "林岱樺"
    $ws2.Cells.Item($r, 6).Value  = $row.Total
    $ws2.Cells.Item($r, 7).Value  = "deposit"
    $ws2.Cells.Item($r, 8).Value  = "normal"

    # Copy the literal date text from sheet1 (J2 = "2012-02-29") so it lands
    # as a plain shared string instead of being auto-parsed into a date
    # serial number.
    $ws1.Range("J2").Copy($ws2.Cells.Item($r, 9))

    $ws2.Cells.Item($r, 10).Value = "林岱樺"
    $ws2.Cells.Item($r, 11).Value = 904
    $ws2.Cells.Item($r, 12).Value = "tmp3bff1"
    $ws2.Cells.Item($r, 13).Value = $row.Idx
    $r = $r + 1
}
